$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.450082182884216
$ws.Range("B1").Value = 1.738661170005798
$ws.Range("C1").Value = 2.331111192703247
$ws.Range("D1").Value = 5.025181293487549
$ws.Range("E1").Value = 1.899722218513489
